# Applies the "fully reorganise client code" version-scheduling edit:
#  - sets completion dates for three tasks (D6, D10, D12)
#  - updates the view state (scroll position + active selection)
#  - resets row 1's explicit height back to the sheet default

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-completed dates (serials 45454 = 2024-06-11, 45456 = 2024-06-13).
# D10/D12 already carry date number formats; D6 needs the plain date format applied,
# matching the other "Completion Date" cells -- reuse D12's existing style so it maps
# onto the same shared cellXf rather than minting a brand-new one.
$ws.Range("D12").Copy()
$ws.Range("D6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D6").Value = 45454
$ws.Range("D10").Value = 45456
$ws.Range("D12").Value = 45454

# Row 1 no longer needs an explicit row height override.
$ws.Rows.Item(1).RowHeight = 14.25

# Update scroll position / current selection to match the saved view.
$ws.Range("B26").Select()
$excel.ActiveWindow.ScrollRow = 7
